$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "priya@msn.usa123"
$ws.Range("B3").Value = "neha008b@12gmail.com"
$ws.Range("B4").Value = "peeyush@yah12oo.india"
$ws.Range("B6").Value = "pallavi.india@ac.indiaa.in"
$ws.Range("B7").Value = "Sampath@hotmail.ind1.us"
$ws.Range("B8").Value = "ramesh.r123g@google.india"
$ws.Range("B5").Value = "manoj_verma@yahoo.g12.india"

$ws.Range("B5").Select()
